$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading paragraph: insert a new run containing just a space " " (sz=28)
#    right before the existing "Product Backlog: " run.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(1)
$headingStart = $headingPara.Range.Start
$headingInsertPoint = $d.Range($headingStart, $headingStart)
$headingInsertPoint.InsertBefore(" ")
# Give the newly inserted run the same explicit sz=28 run formatting as its
# neighbour ("Product Backlog: ") instead of leaving it with no rPr at all.
$headingNewRun = $d.Range($headingStart, $headingStart + 1)
$headingNewRun.Font.Size = 14

# ---------------------------------------------------------------------------
# Helper: build a scoped Range for a given table cell (row/col are 1-based).
# Using $d.Range(start, end) instead of cell.Range.Duplicate() is required so
# that Find.Execute stays scoped to the cell instead of matching elsewhere.
# ---------------------------------------------------------------------------
function Get-CellRange($table, $row, $col) {
    $c = $table.Cell($row, $col)
    return $d.Range($c.Range.Start, $c.Range.End)
}

# XML wrapper used with Range.InsertXML to inject raw WordprocessingML.
function New-XmlPackage($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 2) Row with ID 1 (table row 2): Status "To be started" -> "In Progress"
# ---------------------------------------------------------------------------
$rngR1 = Get-CellRange $t 2 7
$null = $rngR1.Find.Execute("To be started", $true, $false, $false, $false, $false, $true, 1, $false, "In Progress", 1)

# ---------------------------------------------------------------------------
# 3) Row with ID 2 (table row 3): Status "To be started" -> "Done"
# ---------------------------------------------------------------------------
$rngR2 = Get-CellRange $t 3 7
$null = $rngR2.Find.Execute("To be started", $true, $false, $false, $false, $false, $true, 1, $false, "Done", 1)

# ---------------------------------------------------------------------------
# 4) Row with ID 3 (table row 4): Status "To be started" -> "In Progress"
#    and a _GoBack bookmark is added right after the new run.
# ---------------------------------------------------------------------------
$cellR3 = $t.Cell(4, 7)
$rngR3 = $d.Range($cellR3.Range.Start, $cellR3.Range.End)
$rngR3.End = $rngR3.End - 2   # drop trailing paragraph mark + cell mark
$xmlR3 = New-XmlPackage '<w:p w:rsidR="00E8539A" w:rsidRDefault="001A5239" w:rsidP="00CF0139"><w:r><w:t>In Progress</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rngR3.InsertXML($xmlR3)

# ---------------------------------------------------------------------------
# 5) Row with ID 9 (table row 10): the pre-existing _GoBack bookmark that
#    used to sit around the Sprint value "5" is removed (moved away because
#    of the edit above, Word only ever keeps a single _GoBack bookmark).
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(1)
$cellR9 = $t2.Cell(10, 6)
$rngR9 = $cellR9.Range
$rngR9.End = $rngR9.End - 2   # drop trailing paragraph mark + cell mark
$xmlR9 = New-XmlPackage '<w:p w:rsidR="00E8539A" w:rsidRDefault="00F27280" w:rsidP="00CF0139"><w:r><w:t>5</w:t></w:r></w:p>'
$rngR9.InsertXML($xmlR9)

# ---------------------------------------------------------------------------
# 6) Row with ID 10 (table row 11): user story text updated.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(1)
$rngR10c = Get-CellRange $t3 11 3
$null = $rngR10c.Find.Execute("Make a form", $true, $false, $false, $false, $false, $true, 1, $false, "Add category list", 1)

$t4 = $d.Tables.Item(1)
$rngR10d = Get-CellRange $t4 11 4
$null = $rngR10d.Find.Execute("I can store my id and password", $true, $false, $false, $false, $false, $true, 1, $false, "I can add videos ", 1)
